# Ancillary-Structures.docx edit:
#   The broken inline picture (a 1x1 px placeholder image) right under the
#   "Ancillary Structures" heading is replaced by a plain hyperlink whose
#   display text is the image's original URL. The rest of the document
#   (including the existing "here" hyperlink to the NPARKS guidelines
#   further down) is left untouched.

$d = $word.ActiveDocument

$imageUrl = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/BP03_Setbacks_Ancillary_Structures_Substation.jpg?h=100%25&w=100%25"

# The document has exactly one inline picture - the placeholder sitting
# alone in the "FirstParagraph"-styled paragraph right after the heading.
$shape = $d.InlineShapes.Item(1)
$rng = $shape.Range

# Drop the picture and turn the (now empty) range it occupied into a
# hyperlink, using the image URL as both the target address and the
# visible link text.
$shape.Delete()
$d.Hyperlinks.Add($rng, $imageUrl, $null, $null, $imageUrl, $null)
